# Weekly update: insert a new week's worth of Cilantro price records
# (Mercado Mayorista Lo Valledor de Santiago) at the top of the data
# block, pushing the existing rows 1025:1047 down to 1027:1049.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the two new records.
$ws.Rows("1025:1026").Insert()

# New row 1025 - Primera
$ws.Range("A1025").Value = 6
$ws.Range("B1025").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1025").Value = "Metropolitana"
$ws.Range("D1025").Value = 44890
$ws.Range("E1025").Value = 13
$ws.Range("F1025").Value = 100112040
$ws.Range("G1025").Value = "Cilantro"
$ws.Range("H1025").Value = "Sin especificar"
$ws.Range("I1025").Value = "Primera"
$ws.Range("J1025").Value = 570
$ws.Range("K1025").Value = 15000
$ws.Range("L1025").Value = 16000
$ws.Range("M1025").Value = 15614
$ws.Range("N1025").Value = "$/caja 36 atados"
$ws.Range("O1025").Value = "Región Metropolitana"
$ws.Range("P1025").Value = 434
$ws.Range("Q1025").Value = 36
$ws.Range("R1025").Value = "Hortaliza"

# New row 1026 - Segunda
$ws.Range("A1026").Value = 6
$ws.Range("B1026").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1026").Value = "Metropolitana"
$ws.Range("D1026").Value = 44890
$ws.Range("E1026").Value = 13
$ws.Range("F1026").Value = 100112040
$ws.Range("G1026").Value = "Cilantro"
$ws.Range("H1026").Value = "Sin especificar"
$ws.Range("I1026").Value = "Segunda"
$ws.Range("J1026").Value = 150
$ws.Range("K1026").Value = 30000
$ws.Range("L1026").Value = 30000
$ws.Range("M1026").Value = 30000
$ws.Range("N1026").Value = "$/docena de atados"
$ws.Range("O1026").Value = "Región de Valparaíso"
$ws.Range("P1026").Value = 10000
$ws.Range("Q1026").Value = 3
$ws.Range("R1026").Value = "Hortaliza"
